# Applies the three related edits described by the diff:
#  1. The "$ kubectl run ... guestbook --image=ibmcom/guestbook:v1" example command
#     (previously split across many runs, with a "–-generator=run-pod/v1" fragment
#     and a couple of proofErr spell-check wrappers) is collapsed into a single run
#     whose text is "$ kubectl run guestbook --image=ibmcom/guestbook:v1".
#  2. The "$ ibmcloud ks workers mycluster" example command run is split in two:
#     "$ " stays as-is, and "ibmcloud ks workers mycluster" becomes its own run
#     wrapped in the "_GoBack" bookmark.
#  3. Because a document can only have one bookmark with a given name, re-adding
#     "_GoBack" at its new location automatically removes it from its old location
#     (an empty paragraph right before the final section break), which is exactly
#     what the diff shows.

$d = $word.ActiveDocument

# --- Change 1: merge the "kubectl run" command runs into a single run ---
$rng1 = $d.Content
$rng1.Find.Execute("kubectl run", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $rng1.Paragraphs(1).Range
# Exclude the trailing paragraph mark so we only replace the run content.
$target1 = $d.Range($para1.Start, $para1.End - 1)
$target1.Text = "$ kubectl run guestbook --image=ibmcom/guestbook:v1"

# --- Change 2 & 3: split the "ibmcloud ks workers mycluster" run and move the
#     "_GoBack" bookmark onto it (this also removes the stale bookmark at the
#     end of the document since bookmark names are unique). ---
$rng2 = $d.Content
$rng2.Find.Execute("ibmcloud ks workers mycluster", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $rng2)
